$wb = $excel.ActiveWorkbook

# 1. Update status text "Ready for handoff" -> "In Translation" on every sheet
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $v = $cell.Value2
        if ("Ready for handoff" -eq $v) {
            $cell.Value = "In Translation"
        }
    }
}

# 2. Narrow the "Status" / language columns (was 17.2159881591797 chars wide,
#    now 13.4101845877511 chars wide)
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5   # E: zh-cn
$overview.Columns.Item(6).ColumnWidth = 12.5   # F: de-de

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5        # C: Status

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5        # C: Status
